# Add the new survey response row (row 22) submitted by Edward Hira,
# mirroring the columns already populated for earlier participants.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SmartScore" columns for this submission were exported as plain text
# (e.g. "0.561") rather than numbers, so format those cells as Text before
# writing them -- otherwise Excel auto-converts the numeric-looking string
# into a real number and drops the trailing zeros (e.g. "0.730" -> 0.73).
$scoreCells = @("I22", "L22", "O22", "R22", "U22", "X22", "AA22", "AD22", "AG22")
foreach ($ref in $scoreCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("A22").Value = "Edward Hira_20251202_131403"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "Edward Hira"
$ws.Range("D22").Value = 18
$ws.Range("E22").Value = "Male"
$ws.Range("F22").Value = "2025-12-02 13:14:04"
$ws.Range("G22").Value = "{
  ""portion"": 0.8,
  ""diet"": 1.0,
  ""salt"": 0.6,
  ""fat"": 0.8,
  ""natural"": 0.8,
  ""convenience"": 1.0,
  ""price"": 0.8
}"
$ws.Range("H22").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I22").Value = "0.561"
$ws.Range("J22").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K22").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L22").Value = "0.486"
$ws.Range("M22").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("N22").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("O22").Value = "0.473"
$ws.Range("P22").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("Q22").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("R22").Value = "0.595"
$ws.Range("S22").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("T22").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U22").Value = "0.573"
$ws.Range("V22").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("W22").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X22").Value = "0.524"
$ws.Range("Y22").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("Z22").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA22").Value = "0.730"
$ws.Range("AB22").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC22").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD22").Value = "0.600"
$ws.Range("AE22").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AF22").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AG22").Value = "0.573"
$ws.Range("AH22").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# The multi-line JSON in column G makes Excel auto-expand the row height;
# re-run autofit so the new row keeps the sheet's default (non-custom)
# height, matching every other row in this sheet.
$ws.Rows.Item(22).AutoFit()
